$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("BDSBaPCF")

# Set hard coal value to 0 (was 1)
$wsData.Range("B2").Value = 0

# Municipal solid waste (row 17) used to be "=B9" (biomass); now a hardcoded 0
$wsData.Range("B17").Value = 0

# Add explanatory note on the About sheet, row 24, with red font color
$wsAbout.Range("A24").Value = "set value for coal to 0, because HK  has commend control to reduce coal"
$wsAbout.Range("A24").Font.Color = 255

# Scroll / select on the About sheet, then leave BDSBaPCF as the active/selected sheet
# (tabSelected moved from About to BDSBaPCF)
$wsAbout.Activate()
$wsAbout.Range("J25").Select()

$wsData.Activate()
$wsData.Range("J13").Select()

$wb.Save()
